$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $r = $d.Content
    $ok = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Host "NOT FOUND:" $old
    }
}

# Phase 1: replace each original run's text with a unique placeholder marker
Replace-Text ('Contextualizar os fundamentos de Química Orgânica na área de Engenharia de Materiais, para permitir que os estudantes estejam aptos à compreensão das estruturas dos compostos orgânicos e sua influência nas propriedades dos materiais e dos principais mecanismos de reação orgânica a serem utilizados na síntese e processamento de materiais polímeros e outros materiais, como aqueles híbridos, por exemplo.') 'ZMARK1'
Replace-Text ('5840897 - Clodoaldo Saron') 'ZMARK2'
Replace-Text ('1033242 - Fábio Herbst Florenzano') 'ZMARK3'
Replace-Text ('Ligações Químicas e Forças Intermoleculares. Orbitais moleculares e geometria das ligações do carbono. Acidez e basicidade. Reações Orgânicas. Estereoquímica. Principais famílias de compostos de carbono: estrutura e reatividade.') 'ZMARK4'
Replace-Text ('Ligações Químicas Iônicas e Covalentes. Forças intermoleculares e sua relação com as propriedades físicas de compostos orgânicos. Orbitais Moleculares. Geometria das Ligações Covalentes. Conceitos de acidez e basicidade em Química Orgânica. Estereoquímica: diastereoisômeros e enantiômeros. Reações de Substituição Nucleofílica e de Eliminação. Radicais de Carbono e Reatividade. Alcanos, alcenos e alcinos. Composto com grupo acila, álcoois, aminas: formação de poliésteres e poliamidas. Compostos conjugados e aromáticos. Outras famílias de compostos de carbono. Propriedades físicas dos compostos de carbono: relações com a estrutura molecular.') 'ZMARK5'
Replace-Text ('Avaliações envolvendo o conteúdo da disciplina.') 'ZMARK6'
Replace-Text ('Duas avaliações no semestre (P1, P2). MS= (2xP1+3xP2)/5, onde: MS= média do semestre.' + [char]11 + 'MS> ou = 5,0: Aluno Aprovado' + [char]11 + 'MS< 3,0: Aluno Reprovado' + [char]11 + '3,0 < ou = MS < 5,0: Aluno de Recuperação.') 'ZMARK7'
Replace-Text ('Atividade avaliativa versando sobre o conteúdo da disciplina. O aluno será aprovado se apresentar MF (média final) > ou = 5,0. Onde: MF= (MS+PR)/2, onde: MS= média do semestre e PR= prova de recuperação.') 'ZMARK8'
Replace-Text ('McMURRY, J. Química Orgânica. 3ª. Edição. Editora Cengage Learning, 2016.' + [char]11 + '- MORRISON, R.T. e BOYD, R.N. Química Orgânica. 16ª. Edição. Lisboa: Fundacão Calouste Gulbenkian, 2011.' + [char]11 + '- SOLOMONS, T.W.G., FRYHLE, C.B. Química Orgânica 1 e 2. 12ª. Edição, Rio de Janeiro: Gen/LTC Editora, 2018.') 'ZMARK9'

# Phase 2: replace each placeholder marker with the final destination text
Replace-Text 'ZMARK1' ('Ligações Químicas e Forças Intermoleculares. Orbitais moleculares e geometria das ligações do carbono. Acidez e basicidade. Reações Orgânicas. Estereoquímica. Principais famílias de compostos de carbono: estrutura e reatividade.')
Replace-Text 'ZMARK2' ('Contextualizar os fundamentos de Química Orgânica na área de Engenharia de Materiais, para permitir que os estudantes estejam aptos à compreensão das estruturas dos compostos orgânicos e sua influência nas propriedades dos materiais e dos principais mecanismos de reação orgânica a serem utilizados na síntese e processamento de materiais polímeros e outros materiais, como aqueles híbridos, por exemplo.')
Replace-Text 'ZMARK3' ('Ligações Químicas Iônicas e Covalentes. Forças intermoleculares e sua relação com as propriedades físicas de compostos orgânicos. Orbitais Moleculares. Geometria das Ligações Covalentes. Conceitos de acidez e basicidade em Química Orgânica. Estereoquímica: diastereoisômeros e enantiômeros. Reações de Substituição Nucleofílica e de Eliminação. Radicais de Carbono e Reatividade. Alcanos, alcenos e alcinos. Composto com grupo acila, álcoois, aminas: formação de poliésteres e poliamidas. Compostos conjugados e aromáticos. Outras famílias de compostos de carbono. Propriedades físicas dos compostos de carbono: relações com a estrutura molecular.')
Replace-Text 'ZMARK4' ('Avaliações envolvendo o conteúdo da disciplina.')
Replace-Text 'ZMARK5' ('Duas avaliações no semestre (P1, P2). MS= (2xP1+3xP2)/5, onde: MS= média do semestre.' + [char]11 + 'MS> ou = 5,0: Aluno Aprovado' + [char]11 + 'MS< 3,0: Aluno Reprovado' + [char]11 + '3,0 < ou = MS < 5,0: Aluno de Recuperação.')
Replace-Text 'ZMARK6' ('Atividade avaliativa versando sobre o conteúdo da disciplina. O aluno será aprovado se apresentar MF (média final) > ou = 5,0. Onde: MF= (MS+PR)/2, onde: MS= média do semestre e PR= prova de recuperação.')
Replace-Text 'ZMARK7' ('McMURRY, J. Química Orgânica. 3ª. Edição. Editora Cengage Learning, 2016.' + [char]11 + '- MORRISON, R.T. e BOYD, R.N. Química Orgânica. 16ª. Edição. Lisboa: Fundacão Calouste Gulbenkian, 2011.' + [char]11 + '- SOLOMONS, T.W.G., FRYHLE, C.B. Química Orgânica 1 e 2. 12ª. Edição, Rio de Janeiro: Gen/LTC Editora, 2018.')
Replace-Text 'ZMARK8' ('5840897 - Clodoaldo Saron')
Replace-Text 'ZMARK9' ('1033242 - Fábio Herbst Florenzano')

Write-Host "DONE"
